$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "/Large piggy spawner (spawn large pigs)." becomes
#    "+Small piggy spawner (spawn small pigs)." -- the leading tab loses its
#    bold formatting, the bullet glyph changes from "/" to "+", and the rest
#    of the sentence is updated.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(13)
$p1Start = $p1.Range.Start

$tabRange = $d.Range($p1Start, $p1Start + 1)
$tabRange.Font.Bold = 0

$bulletRange = $d.Range($p1Start + 1, $p1Start + 2)
$bulletRange.Text = "+"

$p1b = $d.Paragraphs.Item(13)
$p1b.Range.Find.Execute("Large piggy spawner (spawn large pigs)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Small piggy spawner (spawn small pigs)", 2)

# ---------------------------------------------------------------------------
# 2) The whole "Small piggy spawner (spawn small pigs)." paragraph (the one
#    that used to follow) is removed entirely, and the leading "\t/" of the
#    following "Large piggy (gives 2 points...)" paragraph is removed too --
#    so that paragraph now starts with the old "\t+" that used to belong to
#    the deleted paragraph.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(14)
$bodyStart = $p2.Range.Start + 2          # right after "\t+"
$bodyEnd = $p2.Range.End - 1              # right before the paragraph mark
$d.Range($bodyStart, $bodyEnd).Text = ""

$p2b = $d.Paragraphs.Item(14)
$pilcrowPos = $p2b.Range.End - 1
$d.Range($pilcrowPos, $pilcrowPos + 1).Delete(1, 1)

$p2c = $d.Paragraphs.Item(14)
$prefixStart = $p2c.Range.Start + 2       # right after the surviving "\t+"
$d.Range($prefixStart, $prefixStart + 2).Text = ""

# ---------------------------------------------------------------------------
# 3) Expand the "One shot one kill" remark into the large/small piggy detail.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("One shot one kill", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Small piggies = one shot kill, large piggies = two shot kill", 2)

$d.Save()
